# Hortaliza, Vega Modelo de Temuco - Papa : weekly update
# Inserts 3 new price-report rows at the top of the "Papa" data block
# (rows 712-714), shifting the existing rows 712..809 down to 715..812.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows above the current row 712 (the existing
# data that lived in 712..809 is pushed down to 715..812).
$ws.Rows.Item(712).Insert()
$ws.Rows.Item(713).Insert()
$ws.Rows.Item(714).Insert()

# Row 712 (new)
$ws.Cells.Item(712, 1).Value = 10
$ws.Cells.Item(712, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(712, 3).Value = "La Araucanía"
$ws.Cells.Item(712, 4).Value = 44776
$ws.Cells.Item(712, 5).Value = 9
$ws.Cells.Item(712, 6).Value = 100114001
$ws.Cells.Item(712, 7).Value = "Papa"
$ws.Cells.Item(712, 8).Value = "Rosara"
$ws.Cells.Item(712, 9).Value = "1a (guarda)"
$ws.Cells.Item(712, 10).Value = 400
$ws.Cells.Item(712, 11).Value = 7500
$ws.Cells.Item(712, 12).Value = 7500
$ws.Cells.Item(712, 13).Value = 7500
$ws.Cells.Item(712, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(712, 15).Value = "Provincia de Valdivia"
$ws.Cells.Item(712, 16).Value = 300
$ws.Cells.Item(712, 17).Value = 25
$ws.Cells.Item(712, 18).Value = "Hortaliza"

# Row 713 (new)
$ws.Cells.Item(713, 1).Value = 10
$ws.Cells.Item(713, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(713, 3).Value = "La Araucanía"
$ws.Cells.Item(713, 4).Value = 44776
$ws.Cells.Item(713, 5).Value = 9
$ws.Cells.Item(713, 6).Value = 100114001
$ws.Cells.Item(713, 7).Value = "Papa"
$ws.Cells.Item(713, 8).Value = "Rosara"
$ws.Cells.Item(713, 9).Value = "1a (guarda)"
$ws.Cells.Item(713, 10).Value = 300
$ws.Cells.Item(713, 11).Value = 7000
$ws.Cells.Item(713, 12).Value = 7000
$ws.Cells.Item(713, 13).Value = 7000
$ws.Cells.Item(713, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(713, 15).Value = "Provincia de Valdivia"
$ws.Cells.Item(713, 16).Value = 280
$ws.Cells.Item(713, 17).Value = 25
$ws.Cells.Item(713, 18).Value = "Hortaliza"

# Row 714 (new)
$ws.Cells.Item(714, 1).Value = 10
$ws.Cells.Item(714, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(714, 3).Value = "La Araucanía"
$ws.Cells.Item(714, 4).Value = 44776
$ws.Cells.Item(714, 5).Value = 9
$ws.Cells.Item(714, 6).Value = 100114001
$ws.Cells.Item(714, 7).Value = "Papa"
$ws.Cells.Item(714, 8).Value = "Rosara"
$ws.Cells.Item(714, 9).Value = "2a (guarda)"
$ws.Cells.Item(714, 10).Value = 200
$ws.Cells.Item(714, 11).Value = 4500
$ws.Cells.Item(714, 12).Value = 4500
$ws.Cells.Item(714, 13).Value = 4500
$ws.Cells.Item(714, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(714, 15).Value = "Provincia de Valdivia"
$ws.Cells.Item(714, 16).Value = 180
$ws.Cells.Item(714, 17).Value = 25
$ws.Cells.Item(714, 18).Value = "Hortaliza"
